$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.655.90"
$ws.Range("E2").Value = '  +1.76%  '
$ws.Range("D3").Value = "'2.303.45"
$ws.Range("E3").Value = '  +0.86%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").Value = "'319.52"
$ws.Range("E5").Value = '  +1.52%  '
$ws.Range("D6").Value = "'104.70"
$ws.Range("E6").Value = '  +2.26%  '
$ws.Range("D7").Value = "'0.631"
$ws.Range("E7").Value = '  +0.82%  '
$ws.Range("E8").Value = '  +0.13%  '
$ws.Range("D9").Value = "'0.609"
$ws.Range("E9").Value = '  +1.52%  '
$ws.Range("D10").Value = "'40.14"
$ws.Range("E10").Value = '  +3.76%  '
$ws.Range("D11").Value = "'0.0909"
$ws.Range("E11").Value = '  +0.75%  '
$ws.Range("E12").Value = '  +5.01%  '
$ws.Range("E13").Value = '  +1.20%  '
$ws.Range("D14").Value = "'0.976"
$ws.Range("E14").Value = '  +1.82%  '
$ws.Range("D15").Value = "'15.37"
$ws.Range("E15").Value = '  +1.01%  '
$ws.Range("D16").Value = "'2.653.09"
$ws.Range("E16").Value = '  +1.05%  '
$ws.Range("D17").Value = "'2.298.23"
$ws.Range("E17").Value = '  +0.27%  '
$ws.Range("D18").Value = "'42.765.15"
$ws.Range("E18").Value = '  +2.57%  '
$ws.Range("D19").Value = "'7.53"
$ws.Range("E19").Value = '  +1.64%  '
$ws.Range("D20").Value = "'0.0000105"
$ws.Range("E20").Value = '  +0.77%  '
$ws.Range("D21").Value = "'13.46"
$ws.Range("E21").Value = '  +34.13%  '
$ws.Range("D22").Value = "'73.82"
$ws.Range("E22").Value = '  +1.17%  '
$ws.Range("D23").Value = "'3.57"
$ws.Range("E23").Value = '  -1.00%  '
$ws.Range("D24").Value = "'271.11"
$ws.Range("E24").Value = '  -1.11%  '
$ws.Range("D25").Value = "'2.26"
$ws.Range("E25").Value = '  +0.99%  '
$ws.Range("E26").Value = '  -0.57%  '
$ws.Range("D27").Value = "'10.94"
$ws.Range("E27").Value = '  +2.09%  '
$ws.Range("E28").Value = '  -2.29%  '
$ws.Range("D29").Value = "'22.68"
$ws.Range("E29").Value = '  -1.06%  '
$ws.Range("D30").Value = "'38.08"
$ws.Range("E30").Value = '  +9.60%  '
$ws.Range("D31").Value = "'165.90"
$ws.Range("E31").Value = '  +2.00%  '
$ws.Range("D32").Value = "'6.23"
$ws.Range("E32").Value = '  +7.47%  '
$ws.Range("D33").Value = "'0.0892"
$ws.Range("E33").Value = '  +2.84%  '
$ws.Range("D34").Value = "'0.132"
$ws.Range("E34").Value = '  +0.36%  '
$ws.Range("D35").Value = "'0.115"
$ws.Range("E35").Value = '  +0.64%  '
$ws.Range("E36").Value = '  -11.37%  '
$ws.Range("E37").Value = '  +2.30%  '
$ws.Range("D38").Value = "'0.0355"
$ws.Range("E38").Value = '  +2.84%  '
$ws.Range("E39").Value = '  +2.64%  '
$ws.Range("E40").Value = '  -5.24%  '
$ws.Range("D41").Value = "'1.56"
$ws.Range("E41").Value = '  +7.76%  '
$ws.Range("D42").Value = "'99.61"
$ws.Range("E42").Value = '  -0.39%  '
$ws.Range("D43").Value = "'70.69"
$ws.Range("E43").Value = '  +1.84%  '
$ws.Range("E44").Value = '  +1.31%  '
$ws.Range("D46").Value = "'12.33"
$ws.Range("E46").Value = '  +4.31%  '
$ws.Range("D47").Value = "'82.41"
$ws.Range("E47").Value = '  +9.01%  '
$ws.Range("D48").Value = "'114.88"
$ws.Range("E48").Value = '  -1.04%  '
$ws.Range("E49").Value = '  +1.80%  '
$ws.Range("E50").Value = '  -0.32%  '
$ws.Range("D51").Value = "'1.595.82"
$ws.Range("E51").Value = '  +4.68%  '
